# Apply the "ADDITIONAL SCRAPING" edit described by the diff:
#  - Insert a new "Player Info" sheet at the front with player bio data.
#  - Rename the MATCH_CARD_LINK columns to MATCH_CODE on the existing
#    "ODI Batting" and "ODI Bowling" sheets, and replace the full scorecard
#    URLs with just the numeric match code.
#  - Remove the stray empty INNING_NUMBER cells on "ODI Batting".
#  - Append a new "ODI Batting Extra" sheet at the end with additional
#    per-match batting detail.
#  - Final sheet order: Player Info, ODI Batting, ODI Bowling, ODI Batting Extra.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $text) {
    # Forces the cell to hold a genuine text value (not auto-coerced into a
    # number/percentage/date by Excel's normal "General" type inference),
    # while leaving the cell's number-format/style untouched afterwards.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Existing "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")

Set-TextValue $batting.Cells.Item(1, 4) "MATCH_CODE"

$battingEmptyInningRows = @(2, 5, 8, 9, 11, 12, 13, 20, 23, 29, 35, 41, 43)
foreach ($r in $battingEmptyInningRows) {
    $batting.Cells.Item($r, 2).ClearContents()
}

for ($r = 2; $r -le 43; $r++) {
    $linkCell = $batting.Cells.Item($r, 4)
    $link = $linkCell.Value2
    $code = $link -replace '.*MatchCode=', ''
    Set-TextValue $linkCell $code
}

# ---------------------------------------------------------------------
# 2. Existing "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")

Set-TextValue $bowling.Cells.Item(1, 2) "MATCH_CODE"

for ($r = 2; $r -le 42; $r++) {
    $linkCell = $bowling.Cells.Item($r, 2)
    $link = $linkCell.Value2
    $code = $link -replace '.*MatchCode=', ''
    Set-TextValue $linkCell $code
}

# ---------------------------------------------------------------------
# 3. New "Player Info" sheet, inserted before "ODI Batting"
# ---------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $playerInfoHeaders.Count; $c++) {
    $headerCell = $playerInfo.Cells.Item(1, $c)
    Set-TextValue $headerCell $playerInfoHeaders[$c - 1]
    $headerCell.Font.Bold = $true
    $headerCell.HorizontalAlignment = -4108
    $headerCell.VerticalAlignment = -4160
    $headerCell.Borders.LineStyle = 1
}

$playerInfoRow = @("4380", "Pathira Vasan Dushmantha Chameera", "Right Handed", "Right Arm Fast")
for ($c = 1; $c -le $playerInfoRow.Count; $c++) {
    Set-TextValue $playerInfo.Cells.Item(2, $c) $playerInfoRow[$c - 1]
}

# ---------------------------------------------------------------------
# 4. New "ODI Batting Extra" sheet, appended at the end
# ---------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Add()
$battingExtra.Name = "ODI Batting Extra"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra.Move($null, $lastSheet)

$battingExtraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $battingExtraHeaders.Count; $c++) {
    $headerCell = $battingExtra.Cells.Item(1, $c)
    Set-TextValue $headerCell $battingExtraHeaders[$c - 1]
    $headerCell.Font.Bold = $true
    $headerCell.HorizontalAlignment = -4108
    $headerCell.VerticalAlignment = -4160
    $headerCell.Borders.LineStyle = 1
}

$battingExtraData = @(
    @("4233", "", "", "", "", "NO"),
    @("4449", "9", "1", "0", "3.45%", "NO"),
    @("4450", "", "", "", "", "NO"),
    @("4463", "", "", "", "", "NO"),
    @("4464", "11", "", "", "", "YES"),
    @("4465", "10", "1", "0", "3.78%", "NO"),
    @("4469", "10", "0", "1", "5.81%", "NO"),
    @("4470", "", "", "", "", "NO"),
    @("4471", "10", "1", "1", "4.96%", "NO"),
    @("4480", "9", "0", "0", "0.73%", "NO"),
    @("4482", "10", "", "", "", "NO"),
    @("4485", "9", "0", "0", "", "NO"),
    @("4487", "9", "0", "0", "5.58%", "NO"),
    @("4488", "9", "3", "0", "14.29%", "YES"),
    @("4491", "", "", "", "", "NO"),
    @("4523", "9", "0", "0", "0.39%", "NO"),
    @("4527", "", "", "", "", "NO"),
    @("4594", "10", "1", "0", "3.18%", "NO"),
    @("4597", "", "", "", "", "NO"),
    @("4600", "", "", "", "", "")
)

$r = 2
foreach ($row in $battingExtraData) {
    Set-TextValue $battingExtra.Cells.Item($r, 1) $row[0]

    if ($row[1] -ne "") {
        $battingExtra.Cells.Item($r, 2).Value = [double]$row[1]
    } else {
        $battingExtra.Cells.Item($r, 2).Value = ""
    }

    for ($c = 3; $c -le 6; $c++) {
        Set-TextValue $battingExtra.Cells.Item($r, $c) $row[$c - 1]
    }

    $r = $r + 1
}
